$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Monthly Data Levels")

# Update Service Index (column F) values for rows 242-292
$ws.Cells.Item(242, 6).Value = 83.7
$ws.Cells.Item(243, 6).Value = 84.1
$ws.Cells.Item(244, 6).Value = 93.9
$ws.Cells.Item(245, 6).Value = 92.3
$ws.Cells.Item(246, 6).Value = 90.8
$ws.Cells.Item(247, 6).Value = 104.1
$ws.Cells.Item(248, 6).Value = 105.5
$ws.Cells.Item(249, 6).Value = 99.4
$ws.Cells.Item(250, 6).Value = 108.6
$ws.Cells.Item(251, 6).Value = 110.2
$ws.Cells.Item(252, 6).Value = 110.6
$ws.Cells.Item(253, 6).Value = 116.9
$ws.Cells.Item(254, 6).Value = 98.8
$ws.Cells.Item(255, 6).Value = 103.4
$ws.Cells.Item(256, 6).Value = 111.4
$ws.Cells.Item(257, 6).Value = 101.2
$ws.Cells.Item(258, 6).Value = 107.9
$ws.Cells.Item(259, 6).Value = 112.2
$ws.Cells.Item(260, 6).Value = 108.1
$ws.Cells.Item(261, 6).Value = 109.5
$ws.Cells.Item(262, 6).Value = 118.6
$ws.Cells.Item(263, 6).Value = 106
$ws.Cells.Item(264, 6).Value = 110.9
$ws.Cells.Item(265, 6).Value = 118.8
$ws.Cells.Item(266, 6).Value = 104.2
$ws.Cells.Item(267, 6).Value = 105.9
$ws.Cells.Item(268, 6).Value = 113.5
$ws.Cells.Item(269, 6).Value = 104.7
$ws.Cells.Item(270, 6).Value = 111.9
$ws.Cells.Item(271, 6).Value = 114.2
$ws.Cells.Item(272, 6).Value = 108.7
$ws.Cells.Item(273, 6).Value = 111.6
$ws.Cells.Item(274, 6).Value = 108.7
$ws.Cells.Item(275, 6).Value = 109.5
$ws.Cells.Item(276, 6).Value = 111.6
$ws.Cells.Item(277, 6).Value = 118.4
$ws.Cells.Item(278, 6).Value = 105.1
$ws.Cells.Item(279, 6).Value = 106.1
$ws.Cells.Item(280, 6).Value = 113
$ws.Cells.Item(281, 6).Value = 109.8
$ws.Cells.Item(282, 6).Value = 115.3
$ws.Cells.Item(283, 6).Value = 117.4
$ws.Cells.Item(284, 6).Value = 112.1
$ws.Cells.Item(285, 6).Value = 106.8
$ws.Cells.Item(286, 6).Value = 114.8
$ws.Cells.Item(287, 6).Value = 113.5
$ws.Cells.Item(288, 6).Value = 113.6
$ws.Cells.Item(289, 6).Value = 121.9
$ws.Cells.Item(290, 6).Value = 103.2
$ws.Cells.Item(291, 6).Value = 105.9
$ws.Cells.Item(292, 6).Value = 108.8

# Clear Service Index (column F) values for rows 293-362
$ws.Range("F293:F362").ClearContents()
